# Auto-generated Excel COM-interop script
# Applies updated market price / profit figures to the Leve profit tracking sheets
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 2029.3544
$ws.Range("I15").Value = 2029.3544
$ws.Range("K15").Value = 6088.0632
$ws.Range("M15").Value = -5919.0632
$ws.Range("H80").Value = 354030.03
$ws.Range("I80").Value = 588.8125
$ws.Range("J80").Value = 1296540
$ws.Range("K80").Value = 1766.4375
$ws.Range("L80").Value = 3889620
$ws.Range("M80").Value = -768.4375
$ws.Range("N80").Value = -3891616
$ws.Range("H83").Value = 354030.03
$ws.Range("I83").Value = 588.8125
$ws.Range("J83").Value = 1296540
$ws.Range("K83").Value = 5299.3125
$ws.Range("L83").Value = 11668860
$ws.Range("M83").Value = -307.3125
$ws.Range("N83").Value = -11678844
$ws.Range("H86").Value = 2813.0908
$ws.Range("I86").Value = 1448.2
$ws.Range("J86").Value = 4912.923
$ws.Range("K86").Value = 1448.2
$ws.Range("L86").Value = 4912.923
$ws.Range("M86").Value = -325.2
$ws.Range("N86").Value = -7158.923
$ws.Range("H89").Value = 2813.0908
$ws.Range("I89").Value = 1448.2
$ws.Range("J89").Value = 4912.923
$ws.Range("K89").Value = 7241
$ws.Range("L89").Value = 24564.615
$ws.Range("M89").Value = -1625
$ws.Range("N89").Value = -35796.615
$ws.Range("H121").Value = 867.1905
$ws.Range("J121").Value = 874.2632
$ws.Range("L121").Value = 2622.7896
$ws.Range("N121").Value = -6116.7896
$ws.Range("H135").Value = 2644.5334
$ws.Range("I135").Value = 1894.3077
$ws.Range("J135").Value = 7521
$ws.Range("K135").Value = 17048.7693
$ws.Range("L135").Value = 67689
$ws.Range("M135").Value = -14513.7693
$ws.Range("N135").Value = -72759
$ws.Range("H138").Value = 2714.98
$ws.Range("I138").Value = 1456.5758
$ws.Range("J138").Value = 3334.791
$ws.Range("K138").Value = 4369.7274
$ws.Range("L138").Value = 10004.373
$ws.Range("M138").Value = 770.2726000000002
$ws.Range("N138").Value = -20284.373

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 20002.928
$ws.Range("I32").Value = 17365.824
$ws.Range("K32").Value = 17365.824
$ws.Range("M32").Value = -17078.824
$ws.Range("H61").Value = 1360.7693
$ws.Range("I61").Value = 1192.3572
$ws.Range("J61").Value = 1789.4546
$ws.Range("K61").Value = 1192.3572
$ws.Range("L61").Value = 1789.4546
$ws.Range("M61").Value = -980.3571999999999
$ws.Range("N61").Value = -2213.4546
$ws.Range("H63").Value = 1987.7273
$ws.Range("I63").Value = 1987.7273
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 1987.7273
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -1301.7273
$ws.Range("N63").Value = ""
$ws.Range("H66").Value = 1987.7273
$ws.Range("I66").Value = 1987.7273
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 9938.636500000001
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -6506.636500000001
$ws.Range("N66").Value = ""
$ws.Range("H74").Value = 753.5185
$ws.Range("I74").Value = 760.7778
$ws.Range("J74").Value = 746.2593000000001
$ws.Range("K74").Value = 760.7778
$ws.Range("L74").Value = 746.2593000000001
$ws.Range("M74").Value = 113.2222
$ws.Range("N74").Value = -2494.2593
$ws.Range("H77").Value = 753.5185
$ws.Range("I77").Value = 760.7778
$ws.Range("J77").Value = 746.2593000000001
$ws.Range("K77").Value = 3803.889
$ws.Range("L77").Value = 3731.2965
$ws.Range("M77").Value = 564.1110000000003
$ws.Range("N77").Value = -12467.2965
$ws.Range("H80").Value = 15376.667
$ws.Range("J80").Value = 15376.667
$ws.Range("L80").Value = 15376.667
$ws.Range("N80").Value = -17372.667
$ws.Range("H83").Value = 15376.667
$ws.Range("J83").Value = 15376.667
$ws.Range("L83").Value = 46130.001
$ws.Range("N83").Value = -56114.001
$ws.Range("H97").Value = 474.45456
$ws.Range("I97").Value = 474.45456
$ws.Range("K97").Value = 474.45456
$ws.Range("M97").Value = 21.54543999999999
$ws.Range("H132").Value = 1535.3455
$ws.Range("I132").Value = 1051.3889
$ws.Range("J132").Value = 2452.3157
$ws.Range("K132").Value = 3154.1667
$ws.Range("L132").Value = 7356.9471
$ws.Range("M132").Value = -624.1666999999998
$ws.Range("N132").Value = -12416.9471
$ws.Range("H136").Value = 1360.7693
$ws.Range("I136").Value = 1192.3572
$ws.Range("J136").Value = 1789.4546
$ws.Range("K136").Value = 3577.0716
$ws.Range("L136").Value = 5368.3638
$ws.Range("M136").Value = -1027.0716
$ws.Range("N136").Value = -10468.3638

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H53").Value = 34500
$ws.Range("J53").Value = 34500
$ws.Range("L53").Value = 34500
$ws.Range("N53").Value = -35648
$ws.Range("H59").Value = 60780
$ws.Range("J59").Value = 60780
$ws.Range("L59").Value = 60780
$ws.Range("N59").Value = -62474
$ws.Range("H82").Value = 15651.4
$ws.Range("I82").Value = 11419
$ws.Range("J82").Value = 22000
$ws.Range("K82").Value = 11419
$ws.Range("L82").Value = 22000
$ws.Range("M82").Value = -11036
$ws.Range("N82").Value = -22766
$ws.Range("H85").Value = 15651.4
$ws.Range("I85").Value = 11419
$ws.Range("J85").Value = 22000
$ws.Range("K85").Value = 11419
$ws.Range("L85").Value = 22000
$ws.Range("M85").Value = -10093
$ws.Range("N85").Value = -24652
$ws.Range("H86").Value = 2249.1482
$ws.Range("I86").Value = 2030.8235
$ws.Range("K86").Value = 2030.8235
$ws.Range("M86").Value = -907.8235
$ws.Range("H89").Value = 2249.1482
$ws.Range("I89").Value = 2030.8235
$ws.Range("K89").Value = 10154.1175
$ws.Range("M89").Value = -4538.1175
$ws.Range("H107").Value = 10362.714
$ws.Range("I107").Value = 860.4545000000001
$ws.Range("J107").Value = 45204.332
$ws.Range("K107").Value = 860.4545000000001
$ws.Range("L107").Value = 45204.332
$ws.Range("M107").Value = 1059.5455
$ws.Range("N107").Value = -49044.332

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H59").Value = 3000
$ws.Range("I59").Value = 3000
$ws.Range("K59").Value = 3000
$ws.Range("M59").Value = -1855
$ws.Range("H107").Value = 429
$ws.Range("I107").Value = 1011
$ws.Range("J107").Value = 401.2857
$ws.Range("K107").Value = 1011
$ws.Range("L107").Value = 401.2857
$ws.Range("M107").Value = 909
$ws.Range("N107").Value = -4241.2857
$ws.Range("H132").Value = 1521.2941
$ws.Range("I132").Value = 596.6667
$ws.Range("J132").Value = 3740.4
$ws.Range("K132").Value = 1790.0001
$ws.Range("L132").Value = 11221.2
$ws.Range("M132").Value = 739.9999
$ws.Range("N132").Value = -16281.2

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1209.6875
$ws.Range("I5").Value = 798.0952
$ws.Range("J5").Value = 4090.8333
$ws.Range("K5").Value = 2394.2856
$ws.Range("L5").Value = 12272.4999
$ws.Range("M5").Value = -2282.2856
$ws.Range("N5").Value = -12496.4999
$ws.Range("H68").Value = 931.9293
$ws.Range("I68").Value = 507.84482
$ws.Range("J68").Value = 1531.8536
$ws.Range("K68").Value = 1523.53446
$ws.Range("L68").Value = 4595.560799999999
$ws.Range("M68").Value = -712.5344600000001
$ws.Range("N68").Value = -6217.560799999999
$ws.Range("H71").Value = 931.9293
$ws.Range("I71").Value = 507.84482
$ws.Range("J71").Value = 1531.8536
$ws.Range("K71").Value = 4570.60338
$ws.Range("L71").Value = 13786.6824
$ws.Range("M71").Value = -514.6033800000005
$ws.Range("N71").Value = -21898.6824
$ws.Range("H131").Value = 22220.02
$ws.Range("I131").Value = 83810.836
$ws.Range("J131").Value = 2244.6216
$ws.Range("K131").Value = 251432.508
$ws.Range("L131").Value = 6733.864799999999
$ws.Range("M131").Value = -246392.508
$ws.Range("N131").Value = -16813.8648
$ws.Range("H135").Value = 1209.6875
$ws.Range("I135").Value = 798.0952
$ws.Range("J135").Value = 4090.8333
$ws.Range("K135").Value = 7182.8568
$ws.Range("L135").Value = 36817.4997
$ws.Range("M135").Value = -4647.8568
$ws.Range("N135").Value = -41887.4997

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3828.111
$ws.Range("I80").Value = 3960
$ws.Range("J80").Value = 3168.6667
$ws.Range("K80").Value = 3960
$ws.Range("L80").Value = 3168.6667
$ws.Range("M80").Value = -2962
$ws.Range("N80").Value = -5164.6667
$ws.Range("H83").Value = 3828.111
$ws.Range("I83").Value = 3960
$ws.Range("J83").Value = 3168.6667
$ws.Range("K83").Value = 19800
$ws.Range("L83").Value = 15843.3335
$ws.Range("M83").Value = -14808
$ws.Range("N83").Value = -25827.3335
$ws.Range("H132").Value = 1968.0769
$ws.Range("I132").Value = 1171.5
$ws.Range("J132").Value = 3242.6
$ws.Range("K132").Value = 3514.5
$ws.Range("L132").Value = 9727.799999999999
$ws.Range("M132").Value = -984.5
$ws.Range("N132").Value = -14787.8

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 922.6667
$ws.Range("I61").Value = 922.6667
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 922.6667
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -720.6667
$ws.Range("N61").Value = ""
$ws.Range("H113").Value = 922.6667
$ws.Range("I113").Value = 922.6667
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 922.6667
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 1247.3333
$ws.Range("N113").Value = ""
$ws.Range("H122").Value = 5150.2593
$ws.Range("I122").Value = 6000.6
$ws.Range("J122").Value = 2720.7144
$ws.Range("K122").Value = 18001.8
$ws.Range("L122").Value = 8162.1432
$ws.Range("M122").Value = -15551.8
$ws.Range("N122").Value = -13062.1432
$ws.Range("H132").Value = 2335.6
$ws.Range("I132").Value = 2436.5806
$ws.Range("J132").Value = 2112
$ws.Range("K132").Value = 7309.7418
$ws.Range("L132").Value = 6336
$ws.Range("M132").Value = -4779.7418
$ws.Range("N132").Value = -11396
$ws.Range("H136").Value = 2031.9482
$ws.Range("I136").Value = 1267.8108
$ws.Range("J136").Value = 3378.2856
$ws.Range("K136").Value = 3803.4324
$ws.Range("L136").Value = 10134.8568
$ws.Range("M136").Value = -1253.4324
$ws.Range("N136").Value = -15234.8568

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H104").Value = 38500
$ws.Range("J104").Value = 38500
$ws.Range("L104").Value = 38500
$ws.Range("N104").Value = -45488
$ws.Range("H132").Value = 2095.6667
$ws.Range("I132").Value = 1013.0769
$ws.Range("K132").Value = 3039.2307
$ws.Range("M132").Value = -509.2307000000001
$ws.Range("H136").Value = 2094.2812
$ws.Range("I136").Value = 2373.4443
$ws.Range("J136").Value = 1735.3572
$ws.Range("K136").Value = 7120.3329
$ws.Range("L136").Value = 5206.071599999999
$ws.Range("M136").Value = -4570.3329
$ws.Range("N136").Value = -10306.0716
